$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2,3,4,5,7 keep their own row content but the Id (A), Ost (Q) and
# Nord (R) coordinate columns get reassigned (cyclic permutation of the
# location identifiers among these records).
$ws.Range("A2").Value = 112374124
$ws.Range("Q2").Value = 518088
$ws.Range("R2").Value = 7181677

$ws.Range("A3").Value = 112374157
$ws.Range("Q3").Value = 518111
$ws.Range("R3").Value = 7181672

$ws.Range("A4").Value = 112374285
$ws.Range("Q4").Value = 518160
$ws.Range("R4").Value = 7181550

$ws.Range("A5").Value = 112374261
$ws.Range("Q5").Value = 518145
$ws.Range("R5").Value = 7181574

$ws.Range("A7").Value = 112375655
$ws.Range("Q7").Value = 517882
$ws.Range("R7").Value = 7182353

# Rows 8 and 9 swap their entire record content (species differs between
# them, so every populated cell moves with it).
$ws.Range("A8").Value = 112374196
$ws.Range("B8").Value = 77651
$ws.Range("E8").Value = 230405
$ws.Range("F8").Value = "Garnlav (ssp. sarmentosa)"
$ws.Range("G8").Value = "Alectoria sarmentosa subsp. sarmentosa"
$ws.Range("H8").Value = "(Ach.) Ach."
$ws.Range("J8").Value = ""
$ws.Range("N8").Value = ""
$ws.Range("Q8").Value = 518127
$ws.Range("R8").Value = 7181642
$ws.Range("Z8").ClearContents()
$ws.Range("AB8").ClearContents()
$ws.Range("AF8").Value = ""
$ws.Range("AJ8").ClearContents()
$ws.Range("AK8").ClearContents()
$ws.Range("AM8").ClearContents()
$ws.Range("AO8").ClearContents()

$ws.Range("A9").Value = 112377022
$ws.Range("B9").Value = 89549
$ws.Range("E9").Value = 1108
$ws.Range("F9").Value = "Harticka"
$ws.Range("G9").Value = "Pelloporus leporinus"
$ws.Range("H9").Value = "(Fr.) Krieglst."
$ws.Range("J9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("Q9").Value = 518111
$ws.Range("R9").Value = 7181672
$ws.Range("Z9").Value = "10:43"
$ws.Range("AB9").Value = "10:43"
$ws.Range("AF9").ClearContents()
$ws.Range("AJ9").Value = "gran"
$ws.Range("AK9").Value = "Picea abies"
$ws.Range("AM9").Value = "Stubbe"
$ws.Range("AO9").Value = "Stump # Picea abies"
